$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Company Number column stays text (company numbers can be alphanumeric)
$ws.Range("B2:B21").NumberFormat = "@"

$ws.Range("A2").Value = '4D CAPITAL PROPCO (44) LIMITED'
$ws.Range("B2").Value = '16461269'
$ws.Range("G2").Value = '09:38:51'
$ws.Range("H2").Value = 'Capital'

$ws.Range("A3").Value = 'MYSTIC PARTNERS & CO LTD'
$ws.Range("B3").Value = '16460843'
$ws.Range("G3").Value = '08:52:24'
$ws.Range("H3").Value = 'Partners'

$ws.Range("A4").Value = 'DANCRAFT LLP'
$ws.Range("B4").Value = 'OC456792'
$ws.Range("G4").Value = '10:44:26'
$ws.Range("H4").Value = 'LP'

$ws.Range("A5").Value = 'MOSU VENTURES LTD'
$ws.Range("B5").Value = '16461209'
$ws.Range("G5").Value = '09:38:51'
$ws.Range("H5").Value = 'Ventures'

$ws.Range("A6").Value = 'GLOBAL SQUARE INVESTMENTS LTD'
$ws.Range("B6").Value = '16460615'
$ws.Range("G6").Value = '07:32:31'
$ws.Range("H6").Value = 'Investments'

$ws.Range("A7").Value = 'HAYES PARTNERS LTD'
$ws.Range("B7").Value = '16460412'
$ws.Range("G7").Value = '05:39:32'
$ws.Range("H7").Value = 'Partners'

$ws.Range("A8").Value = 'T GILPIN PHYSIO CONSULTANCY LTD'
$ws.Range("B8").Value = '16460503'
$ws.Range("G8").Value = '06:52:20'
$ws.Range("H8").Value = 'LP'

$ws.Range("A9").Value = 'KVP VENTURES LIMITED'
$ws.Range("B9").Value = '16460424'
$ws.Range("G9").Value = '06:52:20'
$ws.Range("H9").Value = 'Ventures'

$ws.Range("A10").Value = 'SAMVIV PARTNERS LTD'
$ws.Range("B10").Value = '16460672'
$ws.Range("G10").Value = '07:32:31'
$ws.Range("H10").Value = 'Partners'

$ws.Range("A11").Value = 'PULSE SUMMIT CAPITAL LTD'
$ws.Range("B11").Value = '16461206'
$ws.Range("G11").Value = '09:38:51'
$ws.Range("H11").Value = 'Capital'

$ws.Range("A12").Value = 'GS GEPE II SIDECAR IV GP LLP'
$ws.Range("B12").Value = 'SO308186'
$ws.Range("G12").Value = '05:39:32'
$ws.Range("H12").Value = 'GP'

$ws.Range("A13").Value = 'SSMC INVESTMENTS LTD'
$ws.Range("B13").Value = '16461213'
$ws.Range("G13").Value = '09:38:51'
$ws.Range("H13").Value = 'Investments'

$ws.Range("A14").Value = 'FUTURE ENTERPRISES & HOLDING LLP'
$ws.Range("B14").Value = 'OC456787'
$ws.Range("G14").Value = '07:32:31'
$ws.Range("H14").Value = 'LP'

$ws.Range("A15").Value = 'DELWAR INVESTMENTS LIMITED'
$ws.Range("B15").Value = '16460585'
$ws.Range("G15").Value = '07:32:31'
$ws.Range("H15").Value = 'Investments'

$ws.Range("A16").Value = 'DGPI LTD'
$ws.Range("B16").Value = 'SC849118'
$ws.Range("G16").Value = '14:43:22'
$ws.Range("H16").Value = 'GP'

$ws.Range("A17").Value = 'DAVIDSON CAPITAL HOLDINGS LTD'
$ws.Range("B17").Value = 'SC849117'
$ws.Range("G17").Value = '14:43:22'
$ws.Range("H17").Value = 'Capital'

$ws.Range("A18").Value = 'AFROSCOT VENTURES LTD'
$ws.Range("B18").Value = '16462878'
$ws.Range("G18").Value = '14:43:22'
$ws.Range("H18").Value = 'Ventures'

$ws.Range("A19").Value = 'ST GEORGE CAPITAL (LAND) LIMITED'
$ws.Range("B19").Value = '16462880'
$ws.Range("G19").Value = '14:43:22'
$ws.Range("H19").Value = 'Capital'

$ws.Range("A20").Value = 'NEWCO SWANSEA SOCIAL INFRASTRUCTURE LP'
$ws.Range("B20").Value = 'LP024176'
$ws.Range("G20").Value = '21:36:29'
$ws.Range("H20").Value = 'LP'

$ws.Range("A21").Value = 'DCR HOLDINGS & INVESTMENTS LTD'
$ws.Range("B21").Value = '16457272'
$ws.Range("G21").Value = '21:36:29'
$ws.Range("H21").Value = 'Investments'

